# Update automatico via Actualizar@
# Shift the "Fecha" (column D) timestamps down: the newest availability
# check timestamp is written to the top block of rows, and the
# previously-recorded timestamps cascade down into the older blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp   = 44231.76190935122
$shiftedBlock2  = 44230.88231832176
$shiftedBlock3  = 44230.87725172454

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newTimestamp
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $shiftedBlock2
}

for ($r = 30; $r -le 37; $r++) {
    $ws.Cells.Item($r, 4).Value = $shiftedBlock3
}
